$wb = $excel.ActiveWorkbook

$cypherQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN ['Doberman Pinscher'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"
$outputPath = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC17_Canine_Filter_Breed-Doberman_Neo4jData.xlsx"

# --- Add CypherOutput_Message sheet (connection + query + output log) ---
$wsCypherMsg = $wb.Worksheets.Add()
$wsCypherMsg.Name = "CypherOutput_Message"

$wsCypherMsg.Range("A1").Value = "Neo4j_URL:"
$wsCypherMsg.Range("A2").Value = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$wsCypherMsg.Range("A3").Value = "User_name:"
$wsCypherMsg.Range("A4").Value = "neo4j"
$wsCypherMsg.Range("A5").Value = "PWD:"
$wsCypherMsg.Range("A6").Value = "icdcDBneo4j0"
$wsCypherMsg.Range("A7").Value = "Cypher:"
$wsCypherMsg.Range("A8").Value = $cypherQuery
$wsCypherMsg.Range("A9").Value = "Output:"
$wsCypherMsg.Range("A10").Value = $outputPath

$wsCypherMsg.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# --- Add StatOutput sheet (empty) ---
$wsStatOutput = $wb.Worksheets.Add()
$wsStatOutput.Name = "StatOutput"
$wsStatOutput.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# --- Add StatOutput_Message sheet (log: one successful run + one failed/empty-cypher run) ---
$wsStatMsg = $wb.Worksheets.Add()
$wsStatMsg.Name = "StatOutput_Message"

$wsStatMsg.Range("A1").Value = "Neo4j_URL:"
$wsStatMsg.Range("A2").Value = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$wsStatMsg.Range("A3").Value = "User_name:"
$wsStatMsg.Range("A4").Value = "neo4j"
$wsStatMsg.Range("A5").Value = "PWD:"
$wsStatMsg.Range("A6").Value = "icdcDBneo4j0"
$wsStatMsg.Range("A7").Value = "Cypher:"
$wsStatMsg.Range("A8").Value = $cypherQuery
$wsStatMsg.Range("A9").Value = "Output:"
$wsStatMsg.Range("A10").Value = $outputPath
$wsStatMsg.Range("A11").Value = "Cypher query should not be an empty string"
$wsStatMsg.Range("A12").Value = "Neo4j_URL:"
$wsStatMsg.Range("A13").Value = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$wsStatMsg.Range("A14").Value = "User_name:"
$wsStatMsg.Range("A15").Value = "neo4j"
$wsStatMsg.Range("A16").Value = "PWD:"
$wsStatMsg.Range("A17").Value = "icdcDBneo4j0"
$wsStatMsg.Range("A18").Value = "Cypher:"
$wsStatMsg.Range("A19").Value = ""
$wsStatMsg.Range("A20").Value = "Output:"
$wsStatMsg.Range("A21").Value = $outputPath

$wsStatMsg.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
